$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptos price list refresh -- update Price (D) and Volume(1h) (E) columns
# for each row, plus the two swapped rows (46/47: Aptos <-> BabyDogeCoin).
# A handful of Price values are purely decimal-looking strings (e.g. "1.000",
# "0.9987", "24.50") that Excel would otherwise auto-convert to numbers and
# mangle (dropping trailing zeros, switching to scientific notation, etc).
# Force those specific cells to Text format first so the literal string is
# preserved exactly, matching the source data (everything in this sheet is
# stored as text).

$ws.Range("D2").Value = "29.331.20"
$ws.Range("E2").Value = "  -0.25%  "
$ws.Range("D3").Value = "1.842.08"
$ws.Range("E3").Value = "  -0.27%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9987"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "238.95"
$ws.Range("E5").Value = "  -0.85%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6298"
$ws.Range("E6").Value = "  +0.34%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9998"
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07535"
$ws.Range("E8").Value = "  -0.94%  "
$ws.Range("E9").Value = "  -0.82%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.50"
$ws.Range("E10").Value = "  +0.20%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07688"
$ws.Range("E11").Value = "  -0.36%  "
$ws.Range("D12").Value = "1.855.43"
$ws.Range("E12").Value = "  -6.53%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.973"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6778"
$ws.Range("E14").Value = "  -1.40%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.00001037"
$ws.Range("E15").Value = "  +4.67%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "83.07"
$ws.Range("E16").Value = "  +0.19%  "
$ws.Range("D17").Value = "2.116.02"
$ws.Range("E17").Value = "  -6.55%  "
$ws.Range("E18").Value = "  -0.27%  "
$ws.Range("D19").Value = "29.378.56"
$ws.Range("E19").Value = "  -0.94%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "228.67"
$ws.Range("E20").Value = "  -1.25%  "
$ws.Range("E21").Value = "  -0.82%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9995"
$ws.Range("E22").Value = "  -0.12%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.450"
$ws.Range("E23").Value = "  -2.07%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.000"
$ws.Range("E24").Value = "  +0.02%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "156.30"
$ws.Range("E25").Value = "  +1.06%  "
$ws.Range("E26").Value = "  +0.30%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.354"
$ws.Range("E27").Value = "  -1.31%  "
$ws.Range("E28").Value = "  -0.30%  "
$ws.Range("D29").Value = "1.456"
$ws.Range("E29").Value = "  -1.02%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.270"
$ws.Range("E30").Value = "  +1.26%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.05644"
$ws.Range("E31").Value = "  -2.80%  "
$ws.Range("E32").Value = "  -0.19%  "
$ws.Range("E33").Value = "  +0.01%  "
$ws.Range("D34").Value = "1.827"
$ws.Range("E34").Value = "  -1.82%  "
$ws.Range("E35").Value = "  -0.38%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7088"
$ws.Range("E36").Value = "  -1.37%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.590"
$ws.Range("E37").Value = "  -0.38%  "
$ws.Range("D38").Value = "1.240.96"
$ws.Range("E38").Value = "  -0.70%  "
$ws.Range("E39").Value = "  +0.31%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.769"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.236"
$ws.Range("E41").Value = "  +2.35%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9007"
$ws.Range("E42").Value = "  -0.80%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9993"
$ws.Range("E43").Value = "  +0.03%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "101.76"
$ws.Range("E44").Value = "  -0.09%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "65.44"
$ws.Range("E45").Value = "  -2.99%  "
$ws.Range("B46").Value = "Aptos"
$ws.Range("C46").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.103"
$ws.Range("E46").Value = "  -3.03%  "
$ws.Range("B47").Value = "BabyDogeCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000118"
$ws.Range("E47").Value = "  +0.23%  "
$ws.Range("D48").Value = "0.3992"
$ws.Range("E48").Value = "  -0.69%  "
$ws.Range("D49").Value = "1.673"
$ws.Range("E49").Value = "  -1.42%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.920"
$ws.Range("E50").Value = "  -2.64%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.1122"
$ws.Range("E51").Value = "  +0.31%  "
